$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 91 (the "「この宇宙の７つの奇跡」" entry) entirely.
# This shifts all subsequent rows (92..199) up by one, matching the diff.
$ws.Rows.Item(91).Delete()
